$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1333468.4
$ws.Range("I19").Value = 1666772.8
$ws.Range("K19").Value = 1666772.8
$ws.Range("M19").Value = -1666597.8
$ws.Range("H40").Value = 1512.625
$ws.Range("I40").Value = 1575.25
$ws.Range("K40").Value = 1575.25
$ws.Range("M40").Value = -1400.25
$ws.Range("H64").Value = 3090
$ws.Range("H67").Value = 3090
$ws.Range("H76").Value = 3218.4546
$ws.Range("I76").Value = 3040.3
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 3040.3
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -2725.3
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 3218.4546
$ws.Range("I79").Value = 3040.3
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 3040.3
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -1948.3
$ws.Range("N79").Value = -7184
$ws.Range("H132").Value = 21828258
$ws.Range("I132").Value = 23350540
$ws.Range("J132").Value = 8866.666999999999
$ws.Range("K132").Value = 70051620
$ws.Range("L132").Value = 26600.001
$ws.Range("M132").Value = -70049090
$ws.Range("N132").Value = -31660.001
$ws.Range("H133").Value = 31639.523
$ws.Range("J133").Value = 31639.523
$ws.Range("L133").Value = 31639.523
$ws.Range("N133").Value = -41759.523
$ws.Range("H134").Value = 51322.668
$ws.Range("J134").Value = 51322.668
$ws.Range("L134").Value = 51322.668
$ws.Range("N134").Value = -61462.668
$ws.Range("H136").Value = 49281.652
$ws.Range("J136").Value = 49281.652
$ws.Range("L136").Value = 49281.652
$ws.Range("N136").Value = -59481.652
$ws.Range("H139").Value = 42576.844
$ws.Range("J139").Value = 42576.844
$ws.Range("L139").Value = 42576.844
$ws.Range("N139").Value = -52856.844
$ws.Range("H140").Value = 49528.57
$ws.Range("J140").Value = 49528.57
$ws.Range("L140").Value = 49528.57
$ws.Range("N140").Value = -59888.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3154.9285
$ws.Range("I61").Value = 2695.8572
$ws.Range("K61").Value = 2695.8572
$ws.Range("M61").Value = -2483.8572
$ws.Range("H63").Value = 4949770
$ws.Range("I63").Value = 7696942.5
$ws.Range("J63").Value = 4860
$ws.Range("K63").Value = 7696942.5
$ws.Range("L63").Value = 4860
$ws.Range("M63").Value = -7696256.5
$ws.Range("N63").Value = -6232
$ws.Range("H66").Value = 4949770
$ws.Range("I66").Value = 7696942.5
$ws.Range("J66").Value = 4860
$ws.Range("K66").Value = 38484712.5
$ws.Range("L66").Value = 24300
$ws.Range("M66").Value = -38481280.5
$ws.Range("N66").Value = -31164
$ws.Range("H74").Value = 1677.575
$ws.Range("I74").Value = 1261.1212
$ws.Range("K74").Value = 1261.1212
$ws.Range("M74").Value = -387.1212
$ws.Range("H77").Value = 1677.575
$ws.Range("I77").Value = 1261.1212
$ws.Range("K77").Value = 6305.606
$ws.Range("M77").Value = -1937.606
$ws.Range("H88").Value = 13335669
$ws.Range("I88").Value = 16668766
$ws.Range("J88").Value = 3280
$ws.Range("K88").Value = 16668766
$ws.Range("L88").Value = 3280
$ws.Range("M88").Value = -16668360
$ws.Range("N88").Value = -4092
$ws.Range("H91").Value = 13335669
$ws.Range("I91").Value = 16668766
$ws.Range("J91").Value = 3280
$ws.Range("K91").Value = 16668766
$ws.Range("L91").Value = 3280
$ws.Range("M91").Value = -16667362
$ws.Range("N91").Value = -6088
$ws.Range("H136").Value = 3154.9285
$ws.Range("I136").Value = 2695.8572
$ws.Range("K136").Value = 8087.571599999999
$ws.Range("M136").Value = -5537.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1927.64
$ws.Range("I86").Value = 1981.6471
$ws.Range("J86").Value = 1812.875
$ws.Range("K86").Value = 1981.6471
$ws.Range("L86").Value = 1812.875
$ws.Range("M86").Value = -858.6470999999999
$ws.Range("N86").Value = -4058.875
$ws.Range("H89").Value = 1927.64
$ws.Range("I89").Value = 1981.6471
$ws.Range("J89").Value = 1812.875
$ws.Range("K89").Value = 9908.235499999999
$ws.Range("L89").Value = 9064.375
$ws.Range("M89").Value = -4292.235499999999
$ws.Range("N89").Value = -20296.375
$ws.Range("H94").Value = 1228.1923
$ws.Range("I94").Value = 1259.2916
$ws.Range("J94").Value = 855
$ws.Range("K94").Value = 1259.2916
$ws.Range("L94").Value = 855
$ws.Range("M94").Value = -808.2916
$ws.Range("N94").Value = -1757
$ws.Range("H105").Value = 2639.8293
$ws.Range("I105").Value = 2644.4102
$ws.Range("J105").Value = 2550.5
$ws.Range("K105").Value = 2644.4102
$ws.Range("L105").Value = 2550.5
$ws.Range("M105").Value = -897.4101999999998
$ws.Range("N105").Value = -6044.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 71432290
$ws.Range("I62").Value = 71432290
$ws.Range("K62").Value = 71432290
$ws.Range("M62").Value = -71431666
$ws.Range("H65").Value = 71432290
$ws.Range("I65").Value = 71432290
$ws.Range("K65").Value = 357161450
$ws.Range("M65").Value = -357158330
$ws.Range("H122").Value = 2150.9312
$ws.Range("I122").Value = 1750.7826
$ws.Range("J122").Value = 3684.8333
$ws.Range("K122").Value = 5252.3478
$ws.Range("L122").Value = 11054.4999
$ws.Range("M122").Value = -2802.3478
$ws.Range("N122").Value = -15954.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6320.7354
$ws.Range("I70").Value = 6058.6816
$ws.Range("K70").Value = 6058.6816
$ws.Range("M70").Value = -5788.6816
$ws.Range("H73").Value = 6320.7354
$ws.Range("I73").Value = 6058.6816
$ws.Range("K73").Value = 6058.6816
$ws.Range("M73").Value = -5122.6816
$ws.Range("H80").Value = 50002300
$ws.Range("I80").Value = 125001250
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 125001250
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -125000252
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 50002300
$ws.Range("I83").Value = 125001250
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 625006250
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -625001258
$ws.Range("N83").Value = -24984
$ws.Range("H136").Value = 12370.481
$ws.Range("J136").Value = 12370.481
$ws.Range("L136").Value = 37111.443
$ws.Range("N136").Value = -42211.443
